# "Allow row height configuration"
#
# Updates the ToDo list (Sheet1) to reflect that the row-height configuration
# work has been completed:
#   - Three related, now-finished items (rows 5, 6 & 7) are marked "Done",
#     have their manual sort "Order" number cleared, and are hidden from the
#     filtered view (matching how other completed rows are already hidden).
#   - The now-obsolete implementation note on row 6 is removed.
#   - The "Order" numbers for the rows that follow (8, 9 & 10) are shifted
#     down to close the gap left by the rows above.
#   - A stray "Prereq" value on row 17 is cleared.
#   - The open item about replacing the dropdown menu (row 69) is reworded
#     with more detail and its row grows to a 2-line height.
#   - The active cell/selection on the sheet is left on H15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5: "Limit block width to the with of all grid columns" ---
$ws.Range("I5").Clear()
$ws.Range("J5").Value2 = "Done"
$ws.Rows.Item(5).Hidden = $true

# --- Row 6: "Adjust Row height" ---
$ws.Range("I6").Clear()
$ws.Range("J6").Value2 = "Done"
$ws.Range("K6").Clear()
$ws.Rows.Item(6).Hidden = $true

# --- Row 7: 'Fix:  Freeze Header Row makes header row hidden' ---
$ws.Range("I7").Clear()
$ws.Range("J7").Value2 = "Done"
$ws.Rows.Item(7).Hidden = $true

# --- Rows 8-10: close the gap in the "Order" column ---
$ws.Range("I8").Value2 = 2
$ws.Range("I9").Value2 = 3
$ws.Range("I10").Value2 = 4

# --- Row 17: clear the stray "Prereq" value ---
$ws.Range("H17").Clear()

# --- Row 69: reword item and grow the row to fit two lines ---
$ws.Range("B69").Value2 = 'Replace "DropdownMenu" to "Menu" (currently experimental/locked)'
$ws.Rows.Item(69).RowHeight = 30

# --- Leave the sheet's selection on H15 ---
$ws.Activate()
$ws.Range("H15").Select()
